# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1) Update the "datos actualizados" timestamp string (row 1) ---
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 1 de Junio de 2020 a las 07:05"

# --- 2) Update Tailandia (row 81) totals ---
$ws.Cells.Item(81, 2).Value2 = 3082
$ws.Cells.Item(81, 3).Value2 = 1
$ws.Cells.Item(81, 4).Value2 = 2965
$ws.Cells.Item(81, 5).Value2 = 60

# --- 3) Swap whole rows so that country names & their data exchange places ---
# Belice (row 201) <-> Santa Lucia (row 202)
$rowA = $ws.Range("A201:H201").Value2
$rowB = $ws.Range("A202:H202").Value2
$ws.Range("A201:H201").Value2 = $rowB
$ws.Range("A202:H202").Value2 = $rowA

# Seychelles (row 210) <-> Montserrat (row 211)
$rowA = $ws.Range("A210:H210").Value2
$rowB = $ws.Range("A211:H211").Value2
$ws.Range("A210:H210").Value2 = $rowB
$ws.Range("A211:H211").Value2 = $rowA

# San Bartolome (row 215) <-> Bonaire, San Eustaquio y Saba (row 216)
$rowA = $ws.Range("A215:H215").Value2
$rowB = $ws.Range("A216:H216").Value2
$ws.Range("A215:H215").Value2 = $rowB
$ws.Range("A216:H216").Value2 = $rowA
